# Re-apply the table style used by the three data tables on slides 14-16.
# (Google Slides export originally tagged them with the "No Style, Table
# Grid" style {68EAB8BE-5B80-4206-AF5F-D56F063F23EF}; switch them over to
# {D6D0EB7B-D376-454F-98B4-37293C25AE67}.)

$p = $ppt.ActivePresentation

$oldStyleId = "{68EAB8BE-5B80-4206-AF5F-D56F063F23EF}"
$newStyleId = "{D6D0EB7B-D376-454F-98B4-37293C25AE67}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
